$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add headers for the two new columns I (I0) and J (IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting (bold, bordered, centered) used by the other header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Values for the new I0 and IF columns, one entry per data row (rows 2-35)
$iValues = @(2,2,1,3,1,1,1,2,1,1,1,1,1,1,1,1,1,1,1,1,1,1,3,3,1,2,2,1,1,1,1,1,1,1)
$jValues = @(6,7,6,5,2,6,5,6,5,5,7,6,5,5,6,3,7,6,9,5,6,6,6,7,6,7,7,6,5,6,6,5,4,1)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
